$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row -> (Price, Volume(1h)) new values, matching the updated cryptos list.
$updates = @{
    2  = @("43.636.98", "  -0.11%  ")
    3  = @("2.291.45", "  -0.02%  ")
    4  = @($null, "  +0.24%  ")
    5  = @("114.58", "  +18.66%  ")
    6  = @("268.68", "  +0.33%  ")
    7  = @("0.623", "  -0.05%  ")
    8  = @($null, "  +0.05%  ")
    9  = @("0.622", "  +1.62%  ")
    10 = @("48.18", "  +4.44%  ")
    11 = @($null, "  +0.22%  ")
    12 = @("8.81", "  +12.10%  ")
    13 = @("0.106", "  +1.07%  ")
    14 = @("15.61", "  +2.86%  ")
    15 = @("2.634.87", "  +0.03%  ")
    16 = @("0.847", "  -0.35%  ")
    17 = @("2.292.51", "  +0.00%  ")
    18 = @("43.617.07", "  +0.07%  ")
    19 = @($null, "  +1.97%  ")
    20 = @("6.51", "  +4.59%  ")
    21 = @("72.57", "  +0.31%  ")
    22 = @($null, "  -2.07%  ")
    23 = @("232.96", "  -0.09%  ")
    24 = @("9.80", "  +6.55%  ")
    25 = @($null, "  +12.12%  ")
    26 = @($null, "  -0.07%  ")
    27 = @("11.71", "  +4.18%  ")
    28 = @("42.07", "  +3.77%  ")
    29 = @($null, "  -2.14%  ")
    30 = @($null, "  -0.17%  ")
    31 = @("176.53", "  +0.53%  ")
    32 = @("0.0931", "  +4.17%  ")
    33 = @($null, "  -1.35%  ")
    34 = @("5.56", "  +3.46%  ")
    35 = @($null, "  +0.52%  ")
    36 = @($null, "  +9.17%  ")
    37 = @($null, "  +0.55%  ")
    38 = @($null, "  +0.67%  ")
    39 = @("3.80", "  +11.84%  ")
    40 = @("2.44", "  +5.40%  ")
    41 = @("13.89", "  +12.94%  ")
    42 = @($null, "  +2.15%  ")
    43 = @("72.63", "  +11.43%  ")
    44 = @("6.09", "  +16.64%  ")
    45 = @("1.43", "  +5.53%  ")
    46 = @($null, "  +0.04%  ")
    47 = @($null, "  -0.72%  ")
    48 = @("102.80", "  +5.53%  ")
    49 = @($null, "  -1.29%  ")
    50 = @($null, "  +2.80%  ")
    51 = @("0.451", "  +5.09%  ")
}

# Rows whose new Price text is numeric-looking (e.g. "114.58") need the
# cell briefly formatted as Text so Excel keeps the entry as a string
# instead of silently converting it to a number (the source feed always
# stores these as plain text, even when they look numeric). The style is
# then reset back to Normal so no stray number-format/style change is
# left behind on the cell.
$textForceRows = @(5, 6, 7, 9, 10, 12, 13, 14, 16, 20, 21, 23, 24, 27, 28, 31, 32, 34, 39, 40, 41, 43, 44, 45, 48, 51)

foreach ($row in $updates.Keys) {
    $vals = $updates[$row]
    $priceVal = $vals[0]
    $volVal = $vals[1]

    if ($null -ne $priceVal) {
        $priceCell = $ws.Cells.Item($row, 4)
        if ($textForceRows -contains $row) {
            $priceCell.NumberFormat = "@"
            $priceCell.Value = $priceVal
            $priceCell.Style = "Normal"
        }
        else {
            $priceCell.Value = $priceVal
        }
    }
    if ($null -ne $volVal) {
        $ws.Cells.Item($row, 5).Value = $volVal
    }
}

$wb.Save()
